$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 62.85906733333334
$ws.Cells.Item(2, 8).Value = 188.577202
$ws.Cells.Item(2, 9).Value = 0.145580545806332
$ws.Cells.Item(2, 10).Value = 0.145580545806332
$ws.Cells.Item(2, 13).Value = 0.74396
$ws.Cells.Item(2, 14).Value = 2.23188
$ws.Cells.Item(2, 15).Value = 0.006259003216804254
$ws.Cells.Item(2, 16).Value = 0.006259003216804255
$ws.Cells.Item(2, 17).Value = 46.76463173330666
$ws.Cells.Item(2, 18).Value = 420.88168559976
$ws.Cells.Item(2, 19).Value = 0.0009111891045059507
$ws.Cells.Item(2, 20).Value = 0.0009111891045059508

$ws.Cells.Item(3, 7).Value = 62.85906733333334
$ws.Cells.Item(3, 8).Value = 188.577202
$ws.Cells.Item(3, 9).Value = 0.145580545806332
$ws.Cells.Item(3, 10).Value = 0.145580545806332
$ws.Cells.Item(3, 13).Value = 88.14978533333333
$ws.Cells.Item(3, 15).Value = 0.7416121699579786
$ws.Cells.Item(3, 16).Value = 0.7416121699579786
$ws.Cells.Item(3, 17).Value = 5541.013291686879
$ws.Cells.Item(3, 18).Value = 49869.11962518191
$ws.Cells.Item(3, 19).Value = 0.1079643044791007
$ws.Cells.Item(3, 20).Value = 0.1079643044791007

$ws.Cells.Item(4, 7).Value = 62.85906733333334
$ws.Cells.Item(4, 8).Value = 188.577202
$ws.Cells.Item(4, 9).Value = 0.145580545806332
$ws.Cells.Item(4, 10).Value = 0.145580545806332
$ws.Cells.Item(4, 13).Value = 29.76859933333333
$ws.Cells.Item(4, 14).Value = 89.305798
$ws.Cells.Item(4, 15).Value = 0.2504459365921425
$ws.Cells.Item(4, 16).Value = 0.2504459365921425
$ws.Cells.Item(4, 17).Value = 1871.226389913022
$ws.Cells.Item(4, 18).Value = 16841.0375092172
$ws.Cells.Item(4, 19).Value = 0.03646005614406211
$ws.Cells.Item(4, 20).Value = 0.03646005614406211

$ws.Cells.Item(5, 7).Value = 62.85906733333334
$ws.Cells.Item(5, 8).Value = 188.577202
$ws.Cells.Item(5, 9).Value = 0.145580545806332
$ws.Cells.Item(5, 10).Value = 0.145580545806332
$ws.Cells.Item(5, 13).Value = 0.2000323333333334
$ws.Cells.Item(5, 14).Value = 0.6000970000000001
$ws.Cells.Item(5, 15).Value = 0.00168289023307462
$ws.Cells.Item(5, 16).Value = 0.00168289023307462
$ws.Cells.Item(5, 17).Value = 12.57384590984378
$ws.Cells.Item(5, 18).Value = 113.164613188594
$ws.Cells.Item(5, 19).Value = 0.0002449960786631484
$ws.Cells.Item(5, 20).Value = 0.0002449960786631484

$ws.Cells.Item(6, 9).Value = 0.331880415407135
$ws.Cells.Item(6, 10).Value = 0.331880415407135
$ws.Cells.Item(6, 13).Value = 0.74396
$ws.Cells.Item(6, 14).Value = 2.23188
$ws.Cells.Item(6, 15).Value = 0.006259003216804254
$ws.Cells.Item(6, 16).Value = 0.006259003216804255
$ws.Cells.Item(6, 17).Value = 106.60947395168
$ws.Cells.Item(6, 18).Value = 959.48526556512
$ws.Cells.Item(6, 19).Value = 0.00207724058762759
$ws.Cells.Item(6, 20).Value = 0.00207724058762759

$ws.Cells.Item(7, 9).Value = 0.331880415407135
$ws.Cells.Item(7, 10).Value = 0.331880415407135
$ws.Cells.Item(7, 13).Value = 88.14978533333333
$ws.Cells.Item(7, 15).Value = 0.7416121699579786
$ws.Cells.Item(7, 16).Value = 0.7416121699579786
$ws.Cells.Item(7, 18).Value = 113686.7844911845
$ws.Cells.Item(7, 19).Value = 0.2461265550366407
$ws.Cells.Item(7, 20).Value = 0.2461265550366407

$ws.Cells.Item(8, 9).Value = 0.331880415407135
$ws.Cells.Item(8, 10).Value = 0.331880415407135
$ws.Cells.Item(8, 13).Value = 29.76859933333333
$ws.Cells.Item(8, 14).Value = 89.305798
$ws.Cells.Item(8, 15).Value = 0.2504459365921425
$ws.Cells.Item(8, 16).Value = 0.2504459365921425
$ws.Cells.Item(8, 17).Value = 4265.840522615462
$ws.Cells.Item(8, 18).Value = 38392.56470353915
$ws.Cells.Item(8, 19).Value = 0.08311810147322923
$ws.Cells.Item(8, 20).Value = 0.08311810147322923

$ws.Cells.Item(9, 9).Value = 0.331880415407135
$ws.Cells.Item(9, 10).Value = 0.331880415407135
$ws.Cells.Item(9, 13).Value = 0.2000323333333334
$ws.Cells.Item(9, 14).Value = 0.6000970000000001
$ws.Cells.Item(9, 15).Value = 0.00168289023307462
$ws.Cells.Item(9, 16).Value = 0.00168289023307462
$ws.Cells.Item(9, 17).Value = 28.66463496692534
$ws.Cells.Item(9, 18).Value = 257.9817147023281
$ws.Cells.Item(9, 19).Value = 0.0005585183096374151
$ws.Cells.Item(9, 20).Value = 0.0005585183096374152

$ws.Cells.Item(10, 7).Value = 144.7357836666667
$ws.Cells.Item(10, 8).Value = 434.207351
$ws.Cells.Item(10, 9).Value = 0.3352056477733801
$ws.Cells.Item(10, 10).Value = 0.3352056477733801
$ws.Cells.Item(10, 13).Value = 0.74396
$ws.Cells.Item(10, 14).Value = 2.23188
$ws.Cells.Item(10, 15).Value = 0.006259003216804254
$ws.Cells.Item(10, 16).Value = 0.006259003216804255
$ws.Cells.Item(10, 17).Value = 107.6776336166533
$ws.Cells.Item(10, 18).Value = 969.0987025498799
$ws.Cells.Item(10, 19).Value = 0.00209805322770454
$ws.Cells.Item(10, 20).Value = 0.00209805322770454

$ws.Cells.Item(11, 7).Value = 144.7357836666667
$ws.Cells.Item(11, 8).Value = 434.207351
$ws.Cells.Item(11, 9).Value = 0.3352056477733801
$ws.Cells.Item(11, 10).Value = 0.3352056477733801
$ws.Cells.Item(11, 13).Value = 88.14978533333333
$ws.Cells.Item(11, 15).Value = 0.7416121699579786
$ws.Cells.Item(11, 16).Value = 0.7416121699579786
$ws.Cells.Item(11, 17).Value = 12758.42826026844
$ws.Cells.Item(11, 18).Value = 114825.8543424159
$ws.Cells.Item(11, 19).Value = 0.2485925878273863
$ws.Cells.Item(11, 20).Value = 0.2485925878273863

$ws.Cells.Item(12, 7).Value = 144.7357836666667
$ws.Cells.Item(12, 8).Value = 434.207351
$ws.Cells.Item(12, 9).Value = 0.3352056477733801
$ws.Cells.Item(12, 10).Value = 0.3352056477733801
$ws.Cells.Item(12, 13).Value = 29.76859933333333
$ws.Cells.Item(12, 14).Value = 89.305798
$ws.Cells.Item(12, 15).Value = 0.2504459365921425
$ws.Cells.Item(12, 16).Value = 0.2504459365921425
$ws.Cells.Item(12, 17).Value = 4308.58155316901
$ws.Cells.Item(12, 18).Value = 38777.2339785211
$ws.Cells.Item(12, 19).Value = 0.08395089240758
$ws.Cells.Item(12, 20).Value = 0.08395089240758001

$ws.Cells.Item(13, 7).Value = 144.7357836666667
$ws.Cells.Item(13, 8).Value = 434.207351
$ws.Cells.Item(13, 9).Value = 0.3352056477733801
$ws.Cells.Item(13, 10).Value = 0.3352056477733801
$ws.Cells.Item(13, 13).Value = 0.2000323333333334
$ws.Cells.Item(13, 14).Value = 0.6000970000000001
$ws.Cells.Item(13, 15).Value = 0.00168289023307462
$ws.Cells.Item(13, 16).Value = 0.00168289023307462
$ws.Cells.Item(13, 17).Value = 28.95183652367189
$ws.Cells.Item(13, 18).Value = 260.566528713047
$ws.Cells.Item(13, 19).Value = 0.0005641143107092727
$ws.Cells.Item(13, 20).Value = 0.0005641143107092728

$ws.Cells.Item(14, 7).Value = 80.88719666666667
$ws.Cells.Item(14, 8).Value = 242.66159
$ws.Cells.Item(14, 9).Value = 0.187333391013153
$ws.Cells.Item(14, 10).Value = 0.187333391013153
$ws.Cells.Item(14, 13).Value = 0.74396
$ws.Cells.Item(14, 14).Value = 2.23188
$ws.Cells.Item(14, 15).Value = 0.006259003216804254
$ws.Cells.Item(14, 16).Value = 0.006259003216804255
$ws.Cells.Item(14, 17).Value = 60.17683883213333
$ws.Cells.Item(14, 18).Value = 541.5915494891999
$ws.Cells.Item(14, 19).Value = 0.001172520296966174
$ws.Cells.Item(14, 20).Value = 0.001172520296966174

$ws.Cells.Item(15, 7).Value = 80.88719666666667
$ws.Cells.Item(15, 8).Value = 242.66159
$ws.Cells.Item(15, 9).Value = 0.187333391013153
$ws.Cells.Item(15, 10).Value = 0.187333391013153
$ws.Cells.Item(15, 13).Value = 88.14978533333333
$ws.Cells.Item(15, 15).Value = 0.7416121699579786
$ws.Cells.Item(15, 16).Value = 0.7416121699579786
$ws.Cells.Item(15, 17).Value = 7130.189022381782
$ws.Cells.Item(15, 18).Value = 64171.70120143603
$ws.Cells.Item(15, 19).Value = 0.1389287226148509
$ws.Cells.Item(15, 20).Value = 0.1389287226148509

$ws.Cells.Item(16, 7).Value = 80.88719666666667
$ws.Cells.Item(16, 8).Value = 242.66159
$ws.Cells.Item(16, 9).Value = 0.187333391013153
$ws.Cells.Item(16, 10).Value = 0.187333391013153
$ws.Cells.Item(16, 13).Value = 29.76859933333333
$ws.Cells.Item(16, 14).Value = 89.305798
$ws.Cells.Item(16, 15).Value = 0.2504459365921425
$ws.Cells.Item(16, 16).Value = 0.2504459365921425
$ws.Cells.Item(16, 17).Value = 2407.898548766535
$ws.Cells.Item(16, 18).Value = 21671.08693889882
$ws.Cells.Item(16, 19).Value = 0.04691688656727116
$ws.Cells.Item(16, 20).Value = 0.04691688656727116

$ws.Cells.Item(17, 7).Value = 80.88719666666667
$ws.Cells.Item(17, 8).Value = 242.66159
$ws.Cells.Item(17, 9).Value = 0.187333391013153
$ws.Cells.Item(17, 10).Value = 0.187333391013153
$ws.Cells.Item(17, 13).Value = 0.2000323333333334
$ws.Cells.Item(17, 14).Value = 0.6000970000000001
$ws.Cells.Item(17, 15).Value = 0.00168289023307462
$ws.Cells.Item(17, 16).Value = 0.00168289023307462
$ws.Cells.Item(17, 17).Value = 12.57384590984378
$ws.Cells.Item(17, 18).Value = 113.164613188594
$ws.Cells.Item(17, 19).Value = 0.0002449960786631484
$ws.Cells.Item(17, 20).Value = 0.0002449960786631484
